# Update "想去人数" (want-to-go count) values in column F
# on sheets "展览" and "全部类型" to reflect the latest scrape.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 2089
$ws1.Range("F6").Value  = 640
$ws1.Range("F7").Value  = 105
$ws1.Range("F8").Value  = 2078
$ws1.Range("F9").Value  = 10727
$ws1.Range("F15").Value = 7578
$ws1.Range("F17").Value = 723
$ws1.Range("F18").Value = 271
$ws1.Range("F20").Value = 3343

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 2089
$ws4.Range("F6").Value  = 640
$ws4.Range("F8").Value  = 105
$ws4.Range("F9").Value  = 2078
$ws4.Range("F12").Value = 10727
$ws4.Range("F18").Value = 7578
$ws4.Range("F20").Value = 723
$ws4.Range("F21").Value = 271
$ws4.Range("F23").Value = 3343
